$wb = $excel.ActiveWorkbook

# Rename "country_constants" tab to "constants" for consistency
$wsConstants = $wb.Worksheets.Item("country_constants")
$wsConstants.Name = "constants"

$wsDropdown = $wb.Worksheets.Item("dropdown_lists")

# Populate new columns B:D on the dropdown_lists sheet
$wsDropdown.Range("B2").Value = $true
$wsDropdown.Range("C2").Value = "scipy"
$wsDropdown.Range("D2").Value = "None"

$wsDropdown.Range("B3").Value = $false
$wsDropdown.Range("C3").Value = "explicit"

$wsDropdown.Range("C4").Value = "runge_kutta"

# Update selections: constants sheet keeps A23 selected (no longer the active tab)
$wsConstants.Range("A23").Select()

# dropdown_lists becomes the active tab, with F13 selected
$wsDropdown.Range("F13").Select()
